$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-11-17 12:59:02"
$woId = "WO-000002"
$refNum = "#Vearthlyco2037"
$dept = "Shipping"
$customer = "CUS530"
$spdLtl = "SPD"
$notes = "test 1 2 3 residual work order oh FedEx"

# Row 4
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = $woId
$ws.Range("C4").Value = $refNum
$ws.Range("D4").Value = $dept
$ws.Range("E4").Value = $customer
$ws.Range("F4").Value = $spdLtl
$ws.Range("G4").Value = "PPF - Insert"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = $notes
$ws.Range("J4").Value = ""
$ws.Range("J4").Font.Bold = $false

# Row 5
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = $woId
$ws.Range("C5").Value = $refNum
$ws.Range("D5").Value = $dept
$ws.Range("E5").Value = $customer
$ws.Range("F5").Value = $spdLtl
$ws.Range("G5").Value = "FAP - Amazon FBA Product Labeling"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = $notes
$ws.Range("J5").Value = ""
$ws.Range("J5").Font.Bold = $false
